$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the R4_Month / R4_Count header fields to Usage_Date / Usage_Count
$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update selection to match the edited cells
$ws.Range("K1:L1").Select()
